$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The methods-outline paragraph describing snapping-shrimp snaps / diurnal
# sub-sampling / fish-call enumeration is being reorganized:
#   * "To sample for diurnal changes..." becomes its own standalone
#     paragraph (and loses its trailing double space).
#   * "These files were then visually analyzed..." paragraph gains a new
#     sentence after "...to enumerate and identify fish calls" explaining
#     the ACI/SPL low-frequency-band purpose, and the old trailing period
#     is replaced by that new sentence.
#   * The "Snaps produced by snapping shrimp..." sentence moves from being
#     its own (now earlier) paragraph to being appended onto the end of
#     the "These files..." paragraph, right before the bookmark paragraph.
#   * Two new blank paragraphs are added right after the bookmark
#     paragraph.
# ---------------------------------------------------------------------------

# Step 1: split what was one paragraph ("To sample...hours.  These files...
# fish calls.") into two, breaking right before "These files were then".
$splitPoint = $d.Content
$splitPoint.Find.Execute("These files were then", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint.Collapse(1)
$splitPoint.InsertParagraphBefore()

# Step 2: the "To sample for diurnal changes..." paragraph no longer ends
# in two trailing spaces -- just the period.
$trimRng = $d.Content
$trimRng.Find.Execute("2100 hours.  ", $true, $false, $false, $false, $false, $true, 1, $false, "2100 hours.", 2) | Out-Null

# Step 3: insert the new sentence about ACI/SPL right after "...fish calls"
# in place of the old trailing period.
$sentenceRng = $d.Content
$newSentence = " to determine their effect on ACI and SPL within the low frequency band (100 " + [char]0x2013 + " 1000 Hz). "
$sentenceRng.Find.Execute("to enumerate and identify fish calls.", $true, $false, $false, $false, $false, $true, 1, $false, ("to enumerate and identify fish calls" + $newSentence), 2) | Out-Null

# Step 4: move the "Snaps produced by snapping shrimp..." sentence from its
# own paragraph to the end of the "These files..." paragraph (right before
# the new ACI/SPL sentence's following bookmark paragraph).
$snapsSentence = "Snaps produced by snapping shrimp were counted using a detector set to a dB and time threshold determined using Raven Pro 1.5.  "

$snapsRng = $d.Content
$snapsRng.Find.Execute($snapsSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$snapsText = $snapsRng.Text

$insertRng = $d.Content
$insertRng.Find.Execute("low frequency band (100 " + [char]0x2013 + " 1000 Hz). ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertRng.Collapse(0)
$insertRng.InsertAfter($snapsText)

# Step 5: delete the now-relocated sentence from its original paragraph,
# along with that paragraph's trailing mark, so the paragraphs merge away.
$deleteRng = $d.Content
$deleteRng.Find.Execute($snapsSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Range($deleteRng.Start, $deleteRng.End + 1).Delete()

# Step 6: add two new blank paragraphs right after the bookmark paragraph
# (the paragraph immediately following "These files...").
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*These files were then*") {
        $bookmarkPara = $p.Next()
        $bookmarkPara.Range.InsertParagraphAfter()
        $bookmarkPara.Range.InsertParagraphAfter()
        break
    }
}
